$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.999.99"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.295.78"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.56"
$ws.Range("E5").Value = "  -3.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.97"
$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("E7").Value = "  -1.32%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.42"
$ws.Range("E10").Value = "  -2.74%  "

$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.27"
$ws.Range("E12").Value = "  -3.90%  "

$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.986"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.36"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.644.61"
$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.291.24"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.639.48"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.31"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.45"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "267.88"
$ws.Range("E23").Value = "  -2.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -5.34%  "

$ws.Range("E25").Value = "  -2.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("E28").Value = "  +16.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -2.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.37"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.14"
$ws.Range("E31").Value = "  -5.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.69"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0854"
$ws.Range("E33").Value = "  -4.12%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("E34").Value = "  -2.45%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("E35").Value = "  +2.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("E38").Value = "  -2.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.89"
$ws.Range("E41").Value = "  +6.70%  "

$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.67"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.749.40"
$ws.Range("E46").Value = "  +9.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.12"
$ws.Range("E47").Value = "  -3.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.54"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.65"
$ws.Range("E49").Value = "  -7.20%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("E50").Value = "  -3.17%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.15"
$ws.Range("E51").Value = "  -3.39%  "
